$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 5977
$ws.Range("F13").Value = 1607
$ws.Range("F15").Value = 1625
$ws.Range("F16").Value = 563
$ws.Range("F18").Value = 648
$ws.Range("F19").Value = 4636
$ws.Range("F20").Value = 97
$ws.Range("F21").Value = 45
$ws.Range("F23").Value = 3361
$ws.Range("F25").Value = 24
$ws.Range("F27").Value = 13
$ws.Range("F28").Value = 2338
$ws.Range("F30").Value = 341
$ws.Range("F36").Value = 8
$ws.Range("F38").Value = 1268
$ws.Range("F39").Value = 1242

$ws = $wb.Worksheets.Item(2)
$ws.Range("F18").Value = 24
$ws.Range("F21").Value = 241

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 746
$ws.Range("F5").Value = 305

$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 746
$ws.Range("F8").Value = 5977
$ws.Range("F22").Value = 1607
$ws.Range("F24").Value = 1625
$ws.Range("F25").Value = 563
$ws.Range("F27").Value = 648
$ws.Range("F28").Value = 4636
$ws.Range("F30").Value = 3361
$ws.Range("F34").Value = 2338
$ws.Range("F36").Value = 341
$ws.Range("F39").Value = 24
$ws.Range("F41").Value = 241
$ws.Range("F46").Value = 8

